$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: to keep numeric-looking strings (e.g. "1.002") stored as
# plain text (matching the source data which uses dotted display strings,
# not real numbers), we briefly force Text format before assigning the
# value, then reset the cell style back to Normal/General so no stray
# number-format style is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.443.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "

$ws.Range("E7").Value = "  -1.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3389"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07584"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.141"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.37%  "

$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.995"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.959"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.572.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001120"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06739"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.287"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "

$ws.Range("E22").Value = "  -2.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.446.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.334"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.667"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.11"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.004"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.749.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.166"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.975"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.832"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08385"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.375"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02474"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.84%  "

$ws.Range("E39").Value = "  -1.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06525"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.464"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6224"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.01%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.809"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5801"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.073"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.219"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.46%  "

$ws.Range("E51").Value = "  -0.13%  "
